$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the "last updated" timestamp banner (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 22:54"

# ---------------------------------------------------------------------------
# Row 4 - Estados Unidos: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 2763485
$ws.Range("C4").Value = 35632
$ws.Range("D4").Value = 1150498
$ws.Range("E4").Value = 1482421
$ws.Range("G4").Value = 444
$ws.Range("H4").Value = 130566

# ---------------------------------------------------------------------------
# Row 10 - Peru: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 288477
$ws.Range("C10").Value = 3264
$ws.Range("D10").Value = 178245
$ws.Range("E10").Value = 100372
$ws.Range("G10").Value = 183
$ws.Range("H10").Value = 9860

# ---------------------------------------------------------------------------
# Row 17 - Alemania: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = 196296
$ws.Range("C17").Value = 464
$ws.Range("E17").Value = 7437
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 9059

# ---------------------------------------------------------------------------
# Rows 26/27 - Suecia & Egipto swap ranking; write the new country + totals
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "Egipto"
$ws.Range("B26").Value = 69814
$ws.Range("C26").Value = 1503
$ws.Range("D26").Value = 18881
$ws.Range("E26").Value = 47899
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 3034

$ws.Range("A27").Value = "Suecia"
$ws.Range("B27").Value = 69692
$ws.Range("C27").Value = 103
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 5370

# ---------------------------------------------------------------------------
# Rows 50/51 - Armenia & Israel swap ranking; write the new country + totals
# ---------------------------------------------------------------------------
$ws.Range("A50").Value = "Israel"
$ws.Range("B50").Value = 26257
$ws.Range("C50").Value = 1013
$ws.Range("D50").Value = 17452
$ws.Range("E50").Value = 8483
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 322

$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 26065
$ws.Range("C51").Value = 523
$ws.Range("D51").Value = 14563
$ws.Range("E51").Value = 11049
$ws.Range("G51").Value = 10
$ws.Range("H51").Value = 453

# ---------------------------------------------------------------------------
# Row 70 - Costa de Marfil: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B70").Value = 9702
$ws.Range("C70").Value = 203
$ws.Range("D70").Value = 4381
$ws.Range("E70").Value = 5253

# ---------------------------------------------------------------------------
# Row 73 - Uzbekistan: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B73").Value = 8781
$ws.Range("C73").Value = 278
$ws.Range("E73").Value = 2908

# ---------------------------------------------------------------------------
# Row 134 - Ruanda: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B134").Value = 1042
$ws.Range("C134").Value = 17
$ws.Range("D134").Value = 480
$ws.Range("E134").Value = 559
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 3

# ---------------------------------------------------------------------------
# Row 141 - Republica del Chad: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("D141").Value = 785
$ws.Range("E141").Value = 7

# ---------------------------------------------------------------------------
# Row 146 - Santo Tome y Principe: refreshed totals
# ---------------------------------------------------------------------------
$ws.Range("B146").Value = 715
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 252
$ws.Range("E146").Value = 450

# ---------------------------------------------------------------------------
# Rows 164/165 - Namibia & Angola swap ranking; write the new country + totals
# ---------------------------------------------------------------------------
$ws.Range("A164").Value = "Angola"
$ws.Range("B164").Value = 291
$ws.Range("C164").Value = 7
$ws.Range("D164").Value = 97
$ws.Range("E164").Value = 179
$ws.Range("G164").Value = 2
$ws.Range("H164").Value = 15

$ws.Range("A165").Value = "Namibia"
$ws.Range("B165").Value = 285
$ws.Range("C165").Value = 80
$ws.Range("D165").Value = 24
$ws.Range("E165").Value = 261
$ws.Range("H165").Value = 0

# ---------------------------------------------------------------------------
# Rows 203/204 - Laos & Santa Lucia swap ranking (tied totals, label-only swap)
# ---------------------------------------------------------------------------
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"
